$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values can look numeric (e.g. "1.00", "0.553") -- force them to
# remain text so Excel does not coerce them into doubles and mangle the
# formatting / precision. We set NumberFormat to Text, assign the value,
# then restore the default "Normal" style so no stray formatting is left
# behind (matches the source workbook, which has no style on these cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.549.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.581.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.972.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.581.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.884"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.562.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0975"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0807"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.120"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.078.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.832.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.35"
$ws.Range("D50").Style = "Normal"

# Columns B, C, E are unambiguous text (names, URLs, padded percentages)
# and can be assigned directly.
$ws.Range("E2").Value = "  -5.73%  "
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  -2.33%  "
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("E7").Value = "  -4.49%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("E17").Value = "  -5.03%  "
$ws.Range("E18").Value = "  -5.89%  "
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  -5.61%  "
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("E40").Value = "  +5.13%  "
$ws.Range("E41").Value = "  -4.39%  "
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("E43").Value = "  -5.89%  "
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  -6.47%  "
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("E51").Value = "  -5.14%  "
